$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.378.74'
$ws.Range("E2").Value = '  -2.99%  '

$ws.Range("D3").Value = '2.247.30'
$ws.Range("E3").Value = '  -3.74%  '

$ws.Range("E4").Value = '  -0.20%  '

$ws.Range("D5").Value = '233.93'
$ws.Range("E5").Value = '  -1.87%  '

$ws.Range("D6").Value = '0.628'
$ws.Range("E6").Value = '  -5.06%  '

$ws.Range("D7").Value = '69.62'
$ws.Range("E7").Value = '  -3.11%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").Value = '0.560'
$ws.Range("E9").Value = '  -3.77%  '

$ws.Range("D10").Value = '0.0995'
$ws.Range("E10").Value = '  +1.11%  '

$ws.Range("D11").Value = '58.71'
$ws.Range("E11").Value = '  +1.09%  '

$ws.Range("D12").Value = '36.61'
$ws.Range("E12").Value = '  +12.98%  '

$ws.Range("D13").Value = '0.105'
$ws.Range("E13").Value = '  -1.91%  '

$ws.Range("D14").Value = '6.73'
$ws.Range("E14").Value = '  -5.04%  '

$ws.Range("D15").Value = '2.584.86'
$ws.Range("E15").Value = '  -3.60%  '

$ws.Range("D16").Value = '15.09'
$ws.Range("E16").Value = '  -5.84%  '

$ws.Range("D17").Value = '0.854'
$ws.Range("E17").Value = '  -4.10%  '

$ws.Range("D18").Value = '2.253.17'
$ws.Range("E18").Value = '  -4.28%  '

$ws.Range("D19").Value = '42.295.01'
$ws.Range("E19").Value = '  -3.03%  '

$ws.Range("D20").Value = '0.0₃0978'
$ws.Range("E20").Value = '  -2.41%  '

$ws.Range("D21").Value = '6.26'
$ws.Range("E21").Value = '  -4.93%  '

$ws.Range("D22").Value = '73.49'
$ws.Range("E22").Value = '  -5.39%  '

$ws.Range("D23").Value = '234.24'
$ws.Range("E23").Value = '  -6.36%  '

$ws.Range("D24").Value = '2.01'
$ws.Range("E24").Value = '  +5.07%  '

$ws.Range("E25").Value = '  +0.03%  '

$ws.Range("D26").Value = '3.67'
$ws.Range("E26").Value = '  -1.01%  '

$ws.Range("D27").Value = '2.40'
$ws.Range("E27").Value = '  -3.32%  '

$ws.Range("D28").Value = '10.02'
$ws.Range("E28").Value = '  -2.53%  '

$ws.Range("E29").Value = '  -2.04%  '

$ws.Range("D30").Value = '169.87'
$ws.Range("E30").Value = '  -2.84%  '

$ws.Range("D31").Value = '20.59'
$ws.Range("E31").Value = '  -6.60%  '

$ws.Range("D32").Value = '0.122'
$ws.Range("E32").Value = '  -3.96%  '

$ws.Range("E33").Value = '  -5.12%  '

$ws.Range("D34").Value = '0.0729'
$ws.Range("E34").Value = '  -0.12%  '

$ws.Range("D35").Value = '5.37'
$ws.Range("E35").Value = '  +0.58%  '

$ws.Range("D36").Value = '4.71'
$ws.Range("E36").Value = '  -6.39%  '

$ws.Range("D37").Value = '3.63'
$ws.Range("E37").Value = '  -2.56%  '

$ws.Range("D38").Value = '21.87'
$ws.Range("E38").Value = '  +16.80%  '

$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.0277'
$ws.Range("E39").Value = '  +3.60%  '

$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").Value = '2.28'
$ws.Range("E40").Value = '  -3.37%  '

$ws.Range("B41").Value = 'THORChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D41").Value = '6.01'
$ws.Range("E41").Value = '  -5.24%  '

$ws.Range("D42").Value = '65.65'
$ws.Range("E42").Value = '  +1.57%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '9.16'
$ws.Range("E43").Value = '  -0.12%  '

$ws.Range("B44").Value = 'FTXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D44").Value = '4.89'
$ws.Range("E44").Value = '  -11.38%  '

$ws.Range("E45").Value = '  -2.29%  '

$ws.Range("D46").Value = '0.192'
$ws.Range("E46").Value = '  -1.57%  '

$ws.Range("E47").Value = '  +0.04%  '

$ws.Range("E48").Value = '  +12.10%  '

$ws.Range("E49").Value = '  -2.70%  '

$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '2.35'
$ws.Range("E50").Value = '  -3.17%  '

$ws.Range("B51").Value = 'Celestia'
$ws.Range("C51").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D51").Value = '10.09'
$ws.Range("E51").Value = '  +8.65%  '
